# Insert two new data rows at row 413 (pushing existing rows 413..459 down to 415..461)
# and populate them with the new "Alcachofa" price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 413.
$ws.Range("413:414").EntireRow.Insert()

# --- New row 413 ---
$ws.Cells.Item(413, 1).Value = 9
$ws.Cells.Item(413, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(413, 3).Value = "Metropolitana"
$ws.Cells.Item(413, 4).Value = 44783
$ws.Cells.Item(413, 5).Value = 13
$ws.Cells.Item(413, 6).Value = 100112013
$ws.Cells.Item(413, 7).Value = "Alcachofa"
$ws.Cells.Item(413, 8).Value = "Española"
$ws.Cells.Item(413, 9).Value = "Extra"
$ws.Cells.Item(413, 10).Value = 90
$ws.Cells.Item(413, 11).Value = 17000
$ws.Cells.Item(413, 12).Value = 17000
$ws.Cells.Item(413, 13).Value = 17000
$ws.Cells.Item(413, 14).Value = "$/caja 25 unidades"
$ws.Cells.Item(413, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(413, 16).Value = 17000
$ws.Cells.Item(413, 17).Value = 1
$ws.Cells.Item(413, 18).Value = "Hortaliza"

# --- New row 414 ---
$ws.Cells.Item(414, 1).Value = 9
$ws.Cells.Item(414, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(414, 3).Value = "Metropolitana"
$ws.Cells.Item(414, 4).Value = 44783
$ws.Cells.Item(414, 5).Value = 13
$ws.Cells.Item(414, 6).Value = 100112013
$ws.Cells.Item(414, 7).Value = "Alcachofa"
$ws.Cells.Item(414, 8).Value = "Española"
$ws.Cells.Item(414, 9).Value = "Primera"
$ws.Cells.Item(414, 10).Value = 160
$ws.Cells.Item(414, 11).Value = 15000
$ws.Cells.Item(414, 12).Value = 17000
$ws.Cells.Item(414, 13).Value = 16375
$ws.Cells.Item(414, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(414, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(414, 16).Value = 546
$ws.Cells.Item(414, 17).Value = 30
$ws.Cells.Item(414, 18).Value = "Hortaliza"

# Make sure the D column on the two new rows keeps the date number format (s="2")
# used by every other row in this column, by copying the format down from row 412.
$ws.Range("D412").Copy()
$ws.Range("D413:D414").PasteSpecial(-4122)  # xlPasteFormats
